$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 1, shifting all existing data (and the
# header row's formatting) down by 2 rows.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# The insert carried the old header row's style down to row 3. Copy that
# formatting back up onto the new row 1 (which becomes the styled header
# row), then strip row 3 back to the default/unstyled look since it is now
# a plain data row.
$ws.Range("A3:J3").Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)
$ws.Range("A3:J3").ClearFormats()
$excel.CutCopyMode = 0

# New row 1: numeric sequence 0-9 across columns A-J.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9

# New row 2: title row.
$ws.Range("A2").Value = "Holding Screws"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "Steel"
$ws.Range("D2").Value = "Stainless Steel"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
